$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that contain data for rows 11 and 12 (excluding columns whose
# values are identical between the two rows, though swapping all is safe too)
$cols = @("A","B","D","E","F","G","H","Q","R","AC","AX")

foreach ($col in $cols) {
    $addr11 = "$col`11"
    $addr12 = "$col`12"
    $v11 = $ws.Range($addr11).Value2
    $v12 = $ws.Range($addr12).Value2
    $ws.Range($addr11).Value2 = $v12
    $ws.Range($addr12).Value2 = $v11
}
